# Add new way of automation
# Appends 204 new Instagram usernames to column A (rows 229-432),
# growing the shared string table and worksheet data to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUsernames = @(
    'maria.onthemove',
    'hwajess_55',
    'meem.qtr___',
    'its._only',
    'fxwxe',
    'ataag_',
    'vjg1h',
    'mmanoory__',
    '00l0x',
    'samaralmaidany',
    'sh300l11',
    '__g0697',
    'abi_alnasr',
    'mnayer000',
    'phillip_rafco',
    'roula_zuhair',
    '_sairasif_',
    'pallanmark',
    'maha.a.almana',
    'awadallah.hadeel',
    'q_gs_q',
    'premio.lb',
    'malthanii',
    'shamamahhhhhhh',
    'has0osi',
    'nbf1_',
    'whateverj882',
    'deeena_ahmad7',
    'aliqatar9900',
    'burakkilavuz',
    'abeermarmash',
    'xh270',
    'a_b_333',
    'ma_ha7164',
    'yyba89',
    'm86gomaa',
    'mbm_81',
    'yasmine_wohoosh',
    'ajmal.sayyed',
    'a212ali',
    'alkubaisin11',
    'a__1100__',
    'dalal_almosallam',
    'marwa_albukhari',
    'nkammourah',
    'xtbm8',
    'abe93__',
    't1nhn',
    'laadydiary',
    'ha_19933',
    'elena_kanello',
    'fadila.6',
    'nbiiintm',
    'mvneer.a',
    'de.almannaii',
    'delizioso.qa',
    'aeyy66',
    'shop_qatar974',
    'h110',
    'azznaj',
    '91_mba',
    'j_ni',
    'nadejdaay',
    '3bood_dave',
    'memo_senorita',
    'alajmiabdullah422',
    'mashail_al_mtairy9',
    'h_al3bdulla',
    '_.ran10',
    '3li_q3',
    'loyalty_qtr',
    'u0shy',
    'ahmad.m.08',
    'labintm',
    'haneen5184',
    'alhumaidimohammed123',
    '_m6dp',
    'qtr907466',
    'lama5050',
    'ga.1811',
    'almarr2',
    'amera__almarri__',
    'rouwaida_flayfel',
    'saraknio',
    'jx.qtri',
    'alotaiibii2_',
    'izzzd_',
    'mmn3060',
    'imixq',
    'amooona1983',
    'sandra1462015',
    '__soso_00',
    '09sep_',
    'blue_lagoon_spa',
    'rkr.10',
    'afa.althanii',
    'x.m.91.x',
    'qtrimem',
    'lattiqat',
    'mbiomy',
    'aishaa.9x',
    'soso_events',
    'phosphenes_qatar',
    'faisal_almehshadi',
    'shintojoseph777',
    'lmeybs',
    'mba.888',
    's71y',
    'm__6__1',
    '_rawanadil',
    'j9xd',
    'hamasaki_1990',
    '24.11.96_',
    'nss__1991',
    'wyi2l',
    'fo0oz2012',
    'lolo4033',
    'ghalya.alobaidly',
    'burger_showqtr',
    'mralfaisall',
    'mm.m812',
    'ercanaykan7',
    '_i3bj',
    'meryemguaid2019',
    '_qxbl',
    'glassjaronsale',
    'makeupabeer_qtr9',
    'aiiishaaa5',
    'ads.qatar',
    'totah_75_',
    'ewan_alhajri2',
    'x.0096',
    '0b6__',
    'mbxre',
    'shmookal3az',
    'lozan.rr',
    'e1002xx',
    'etoo87',
    'tmald_',
    '5.6f3',
    'maryamkhan7629',
    'jiggsylog',
    'tathkaralhamad_',
    'kelvincee911',
    'noura_albuenain',
    'almahaalkuwari94',
    'saraahkhoja',
    'bo_mohammad2121',
    '_90h7',
    'fooz___qtr',
    'matra.alahmed',
    'n_s_qtr11',
    'maialkuwari',
    'afshan1010',
    'healthybody_sonia4',
    'jurida_doko',
    'pergosa.tr',
    'ma95__q',
    'hustohqatar',
    '8llhy',
    'amelcommunication007',
    'joyaksjoyaks',
    '1mqr',
    'chelseabridge.qa',
    'znsour',
    'masviyy',
    'arfaoui_hajer',
    'wmatouk',
    'mohannadnaim',
    'qata_202020',
    'mojarrdi.81',
    'aiijx21',
    '_maalmaadeed',
    'ixxsy',
    'ialmari_',
    'elitebreezeltd',
    'nvxk_',
    'fbintalb9',
    'shan.sampath',
    'i7uii_',
    '_budouralturki',
    'j.vxyy',
    'se_1001',
    'sama.qtr87',
    'alreem930124003',
    'umslman8227',
    'vlez18',
    'hamada_cherkaoui_sellami',
    'iix.aey01',
    'kaak23kaak',
    'ab0902_',
    'rix_o0',
    'salq88__',
    'jabal7182',
    'khulood85',
    'ifoay',
    'ikram_guerfala',
    'photo_lovers9',
    'aishalthani',
    'q_7_q_',
    'novalmarrri',
    'singhshivendra18',
    'saratmahmoud',
    'saracabral27'
)

$startRow = 229
for ($i = 0; $i -lt $newUsernames.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newUsernames[$i]
}

$endRow = $startRow + $newUsernames.Count - 1

# Update the view/selection to mirror where Excel would land after
# pasting this block of new rows.
$selRange = "A" + $startRow + ":A" + $endRow
$ws.Range($selRange).Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 357
$win.ScrollColumn = 1
